$wb = $excel.ActiveWorkbook

# ===== Sheet: LP1912 =====
$ws = $wb.Worksheets.Item("LP1912")
$ws.Range("A2").Value = "Última actualización: 10:26:25"
$ws.Range("A3").Value = "Total filas: 144"
$ws.Range("C49").Value = "17X38_ROMERO"
$ws.Range("C50").Value = "27_EL RETIRO"
$ws.Range("C69").Value = "215A_EL PATO"
$ws.Range("C70").Value = "14_ABASTO"
$ws.Range("A71").Value = "08:30:59"
$ws.Range("C71").Value = "215C_EL PATO"
$ws.Range("D71").Value = 3
$ws.Range("A72").Value = "08:00:50"
$ws.Range("C72").Value = "23_HERNANDEZ"
$ws.Range("D72").Value = 33
$ws.Range("A73").Value = "08:30:59"
$ws.Range("C73").Value = "23_HERNANDEZ"
$ws.Range("D73").Value = 4
$ws.Range("A74").Value = "07:48:14"
$ws.Range("C74").Value = "215C_EL PATO"
$ws.Range("D74").Value = 46
$ws.Range("C85").Value = "215A_EL PATO"
$ws.Range("C86").Value = "215B_EL PATO"
$ws.Range("A87").Value = "08:48:29"
$ws.Range("C87").Value = "215B_EL PATO"
$ws.Range("D87").Value = 11
$ws.Range("A88").Value = "08:30:59"
$ws.Range("C88").Value = "16_P MOR-SANTA ANA"
$ws.Range("D88").Value = 29
$ws.Range("A89").Value = "08:30:59"
$ws.Range("C89").Value = "16_SANTA ANA"
$ws.Range("D89").Value = 30
$ws.Range("A90").Value = "07:48:14"
$ws.Range("C90").Value = "215B_EL PATO"
$ws.Range("D90").Value = 72
$ws.Range("C94").Value = "17X38_ROMERO"
$ws.Range("C95").Value = "23_HERNANDEZ"
$ws.Range("A96").Value = "08:00:50"
$ws.Range("C96").Value = "27_EL RETIRO"
$ws.Range("D96").Value = 74
$ws.Range("A97").Value = "08:48:29"
$ws.Range("C97").Value = "11_ETCHEVERRY"
$ws.Range("D97").Value = 26
$ws.Range("A106").Value = "08:30:59"
$ws.Range("C106").Value = "15_ABASTO"
$ws.Range("D106").Value = 63
$ws.Range("A107").Value = "08:56:14"
$ws.Range("C107").Value = "23_HERNANDEZ"
$ws.Range("D107").Value = 37
$ws.Range("A118").Value = "09:31:15"
$ws.Range("C118").Value = "23_HERNANDEZ"
$ws.Range("D118").Value = 32
$ws.Range("A119").Value = "08:30:59"
$ws.Range("C119").Value = "215C_EL PATO"
$ws.Range("D119").Value = 93
$ws.Range("A122").Value = "08:30:59"
$ws.Range("C122").Value = "10_OLMOS"
$ws.Range("D122").Value = 109
$ws.Range("A123").Value = "09:31:15"
$ws.Range("C123").Value = "17_ROMERO"
$ws.Range("D123").Value = 48
$ws.Range("A125").Value = "10:26:25"
$ws.Range("B125").Value = "10:32"
$ws.Range("D125").Value = 6
$ws.Range("A126").Value = "08:56:14"
$ws.Range("B126").Value = "10:33"
$ws.Range("D126").Value = 97
$ws.Range("A127").Value = "10:26:25"
$ws.Range("B127").Value = "10:33"
$ws.Range("C127").Value = "23_HERNANDEZ"
$ws.Range("D127").Value = 7
$ws.Range("A128").Value = "09:31:15"
$ws.Range("B128").Value = "10:34"
$ws.Range("D128").Value = 63
$ws.Range("A129").Value = "10:26:25"
$ws.Range("B129").Value = "10:34"
$ws.Range("C129").Value = "15_ABASTO"
$ws.Range("D129").Value = 8
$ws.Range("A130").Value = "08:48:29"
$ws.Range("B130").Value = "10:36"
$ws.Range("C130").Value = "14_ABASTO"
$ws.Range("D130").Value = 108
$ws.Range("A131").Value = "10:26:25"
$ws.Range("B131").Value = "10:41"
$ws.Range("C131").Value = "16_SANTA ANA"
$ws.Range("D131").Value = 15
$ws.Range("A132").Value = "10:26:25"
$ws.Range("B132").Value = "10:44"
$ws.Range("C132").Value = "10_OLMOS"
$ws.Range("D132").Value = 18
$ws.Range("A133").Value = "10:26:25"
$ws.Range("B133").Value = "10:49"
$ws.Range("C133").Value = "15_ABASTO"
$ws.Range("D133").Value = 23
$ws.Range("A134").Value = "10:26:25"
$ws.Range("B134").Value = "10:51"
$ws.Range("C134").Value = "16_P MOR-SANTA ANA"
$ws.Range("D134").Value = 25
$ws.Range("A135").Value = "10:26:25"
$ws.Range("B135").Value = "10:56"
$ws.Range("C135").Value = "14_ABASTO"
$ws.Range("D135").Value = 30
$ws.Range("A136").Value = "10:26:25"
$ws.Range("B136").Value = "10:57"
$ws.Range("C136").Value = "27_EL RETIRO"
$ws.Range("D136").Value = 31
$ws.Range("E136").Value = "LP1912"
$ws.Range("A137").Value = "10:26:25"
$ws.Range("B137").Value = "11:01"
$ws.Range("C137").Value = "16_SANTA ANA"
$ws.Range("D137").Value = 35
$ws.Range("E137").Value = "LP1912"
$ws.Range("A138").Value = "10:26:25"
$ws.Range("B138").Value = "11:03"
$ws.Range("C138").Value = "23_HERNANDEZ"
$ws.Range("D138").Value = 37
$ws.Range("E138").Value = "LP1912"
$ws.Range("A139").Value = "10:26:25"
$ws.Range("B139").Value = "11:04"
$ws.Range("C139").Value = "17_ROMERO"
$ws.Range("D139").Value = 38
$ws.Range("E139").Value = "LP1912"
$ws.Range("A140").Value = "10:26:25"
$ws.Range("B140").Value = "11:08"
$ws.Range("C140").Value = "225_C ROCA-H SUR"
$ws.Range("D140").Value = 42
$ws.Range("E140").Value = "LP1912"
$ws.Range("A141").Value = "10:26:25"
$ws.Range("B141").Value = "11:19"
$ws.Range("C141").Value = "215C_EL PATO"
$ws.Range("D141").Value = 53
$ws.Range("E141").Value = "LP1912"
$ws.Range("A142").Value = "10:26:25"
$ws.Range("B142").Value = "11:20"
$ws.Range("C142").Value = "11_ETCHEVERRY"
$ws.Range("D142").Value = 54
$ws.Range("E142").Value = "LP1912"
$ws.Range("A143").Value = "09:31:15"
$ws.Range("B143").Value = "11:21"
$ws.Range("C143").Value = "11_ETCHEVERRY"
$ws.Range("D143").Value = 110
$ws.Range("E143").Value = "LP1912"
$ws.Range("A144").Value = "10:26:25"
$ws.Range("B144").Value = "11:33"
$ws.Range("C144").Value = "215A_EL PATO"
$ws.Range("D144").Value = 67
$ws.Range("E144").Value = "LP1912"
$ws.Range("A145").Value = "10:26:25"
$ws.Range("B145").Value = "11:44"
$ws.Range("C145").Value = "215B_EL PATO"
$ws.Range("D145").Value = 78
$ws.Range("E145").Value = "LP1912"
$ws.Range("A146").Value = "10:26:25"
$ws.Range("B146").Value = "11:51"
$ws.Range("C146").Value = "16_P MOR-SANTA ANA"
$ws.Range("D146").Value = 85
$ws.Range("E146").Value = "LP1912"
$ws.Range("A147").Value = "10:26:25"
$ws.Range("B147").Value = "11:56"
$ws.Range("C147").Value = "225_GOMEZ"
$ws.Range("D147").Value = 90
$ws.Range("E147").Value = "LP1912"
$ws.Range("A148").Value = "10:26:25"
$ws.Range("B148").Value = "12:04"
$ws.Range("C148").Value = "17_ROMERO"
$ws.Range("D148").Value = 98
$ws.Range("E148").Value = "LP1912"
$ws.Range("A149").Value = "10:26:25"
$ws.Range("B149").Value = "12:20"
$ws.Range("C149").Value = "10_OLMOS"
$ws.Range("D149").Value = 114
$ws.Range("E149").Value = "LP1912"

# ===== Sheet: LP1912-215 =====
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Range("A2").Value = "Última actualización: 10:26:25"
$ws.Range("A3").Value = "Total filas: 28"
$ws.Range("C23").Value = "215B_EL PATO"
$ws.Range("C24").Value = "215A_EL PATO"
$ws.Range("A31").Value = "10:26:25"
$ws.Range("D31").Value = 53
$ws.Range("A32").Value = "10:26:25"
$ws.Range("B32").Value = "11:33"
$ws.Range("C32").Value = "215A_EL PATO"
$ws.Range("D32").Value = 67
$ws.Range("E32").Value = "LP1912"
$ws.Range("A33").Value = "10:26:25"
$ws.Range("B33").Value = "11:44"
$ws.Range("C33").Value = "215B_EL PATO"
$ws.Range("D33").Value = 78
$ws.Range("E33").Value = "LP1912"

# ===== Sheet: 6203-6173 =====
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Range("A2").Value = "Última actualización: 10:26:25"
$ws.Range("A3").Value = "Total filas: 25"
$ws.Range("A30").Value = "10:26:25"
$ws.Range("B30").Value = "11:56"
$ws.Range("C30").Value = "215C_LA PLATA"
$ws.Range("D30").Value = 90
$ws.Range("E30").Value = "L6203"
